$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 172.90909
$ws.Range("I5").Value = 207.75
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 207.75
$ws.Range("L5").Value = 80
$ws.Range("M5").Value = -92.75
$ws.Range("N5").Value = -310

$ws.Range("H18").Value = 1715.25
$ws.Range("I18").Value = 1787
$ws.Range("K18").Value = 1787
$ws.Range("M18").Value = -1503

$ws.Range("H33").Value = 58373.57
$ws.Range("I33").Value = 1364.8
$ws.Range("K33").Value = 1364.8
$ws.Range("M33").Value = -1135.8

$ws.Range("H43").Value = 1148.6666
$ws.Range("I43").Value = 1235.125
$ws.Range("J43").Value = 975.75
$ws.Range("K43").Value = 1235.125
$ws.Range("L43").Value = 975.75
$ws.Range("M43").Value = -1166.125
$ws.Range("N43").Value = -1113.75

$ws.Range("H51").Value = 464579.62
$ws.Range("I51").Value = 557045.5600000001
$ws.Range("K51").Value = 557045.5600000001
$ws.Range("M51").Value = -556561.5600000001

$ws.Range("H62").Value = 1075.1666
$ws.Range("I62").Value = 740.8182
$ws.Range("J62").Value = 1600.5714
$ws.Range("K62").Value = 740.8182
$ws.Range("L62").Value = 1600.5714
$ws.Range("M62").Value = -116.8182
$ws.Range("N62").Value = -2848.5714

$ws.Range("H65").Value = 1075.1666
$ws.Range("I65").Value = 740.8182
$ws.Range("J65").Value = 1600.5714
$ws.Range("K65").Value = 3704.091
$ws.Range("L65").Value = 8002.857
$ws.Range("M65").Value = -584.0910000000003
$ws.Range("N65").Value = -14242.857

$ws.Range("H107").Value = 823.3077
$ws.Range("I107").Value = 671.4
$ws.Range("J107").Value = 1329.6666
$ws.Range("K107").Value = 671.4
$ws.Range("L107").Value = 1329.6666
$ws.Range("M107").Value = 1248.6
$ws.Range("N107").Value = -5169.6666

$ws.Range("H116").Value = 2105172.5
$ws.Range("I116").Value = 35715684
$ws.Range("J116").Value = 4515.625
$ws.Range("K116").Value = 35715684
$ws.Range("L116").Value = 4515.625
$ws.Range("M116").Value = -35712242
$ws.Range("N116").Value = -11399.625

$ws.Range("H132").Value = 4632883
$ws.Range("I132").Value = 11114211
$ws.Range("K132").Value = 33342633
$ws.Range("M132").Value = -33340103

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 959.56525
$ws.Range("I2").Value = 597.4375
$ws.Range("J2").Value = 1787.2858
$ws.Range("K2").Value = 597.4375
$ws.Range("L2").Value = 1787.2858
$ws.Range("M2").Value = -484.4375
$ws.Range("N2").Value = -2013.2858

$ws.Range("H102").Value = 1667.3334
$ws.Range("I102").Value = 1679.2858
$ws.Range("K102").Value = 1679.2858
$ws.Range("M102").Value = -57.28580000000011

$ws.Range("H116").Value = 959.56525
$ws.Range("I116").Value = 597.4375
$ws.Range("J116").Value = 1787.2858
$ws.Range("K116").Value = 597.4375
$ws.Range("L116").Value = 1787.2858
$ws.Range("M116").Value = 1696.5625
$ws.Range("N116").Value = -6375.2858

$ws.Range("H132").Value = 1024.8837
$ws.Range("I132").Value = 632.2564
$ws.Range("K132").Value = 1896.7692
$ws.Range("M132").Value = 633.2308

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 959.56525
$ws.Range("I3").Value = 597.4375
$ws.Range("J3").Value = 1787.2858
$ws.Range("K3").Value = 597.4375
$ws.Range("L3").Value = 1787.2858
$ws.Range("M3").Value = -483.4375
$ws.Range("N3").Value = -2015.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2092
$ws.Range("I58").Value = 1366.5
$ws.Range("J58").Value = 2636.125
$ws.Range("K58").Value = 1366.5
$ws.Range("L58").Value = 2636.125
$ws.Range("M58").Value = -1163.5
$ws.Range("N58").Value = -3042.125

$ws.Range("H107").Value = 795.2917
$ws.Range("I107").Value = 333.16666
$ws.Range("J107").Value = 2181.6667
$ws.Range("K107").Value = 333.16666
$ws.Range("L107").Value = 2181.6667
$ws.Range("M107").Value = 1586.83334
$ws.Range("N107").Value = -6021.6667

$ws.Range("H132").Value = 3117.8
$ws.Range("I132").Value = 1197
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 3591
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -1061
$ws.Range("N132").Value = -23057

$ws.Range("H134").Value = 2233.4285
$ws.Range("I134").Value = 1812.8572
$ws.Range("J134").Value = 2864.2856
$ws.Range("K134").Value = 5438.571599999999
$ws.Range("L134").Value = 8592.856800000001
$ws.Range("M134").Value = -2903.571599999999
$ws.Range("N134").Value = -13662.8568

$ws.Range("H136").Value = 2092
$ws.Range("I136").Value = 1366.5
$ws.Range("J136").Value = 2636.125
$ws.Range("K136").Value = 4099.5
$ws.Range("L136").Value = 7908.375
$ws.Range("M136").Value = -1549.5
$ws.Range("N136").Value = -13008.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4333.5
$ws.Range("I68").Value = 999
$ws.Range("J68").Value = 5000.4
$ws.Range("K68").Value = 2997
$ws.Range("L68").Value = 15001.2
$ws.Range("M68").Value = -2186
$ws.Range("N68").Value = -16623.2

$ws.Range("H71").Value = 4333.5
$ws.Range("I71").Value = 999
$ws.Range("J71").Value = 5000.4
$ws.Range("K71").Value = 8991
$ws.Range("L71").Value = 45003.6
$ws.Range("M71").Value = -4935
$ws.Range("N71").Value = -53115.6

$ws.Range("H113").Value = 410.52
$ws.Range("I113").Value = 355.29413
$ws.Range("J113").Value = 438.9697
$ws.Range("K113").Value = 1065.88239
$ws.Range("L113").Value = 1316.9091
$ws.Range("M113").Value = 1104.11761
$ws.Range("N113").Value = -5656.9091

$ws.Range("H131").Value = 945.5961
$ws.Range("I131").Value = 493.33334
$ws.Range("J131").Value = 1004.587
$ws.Range("K131").Value = 1480.00002
$ws.Range("L131").Value = 3013.761
$ws.Range("M131").Value = 3559.99998
$ws.Range("N131").Value = -13093.761

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5859.9585
$ws.Range("I113").Value = 9123.076999999999
$ws.Range("J113").Value = 2003.5454
$ws.Range("K113").Value = 9123.076999999999
$ws.Range("L113").Value = 2003.5454
$ws.Range("M113").Value = -6953.076999999999
$ws.Range("N113").Value = -6343.5454

$ws.Range("H132").Value = 2433.4814
$ws.Range("I132").Value = 2141.3
$ws.Range("J132").Value = 3268.2856
$ws.Range("K132").Value = 6423.900000000001
$ws.Range("L132").Value = 9804.856800000001
$ws.Range("M132").Value = -3893.900000000001
$ws.Range("N132").Value = -14864.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 474.3913
$ws.Range("I22").Value = 469
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 469
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -174
$ws.Range("N22").Value = -1090

$ws.Range("H27").Value = 474.3913
$ws.Range("I27").Value = 469
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 469
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -362
$ws.Range("N27").Value = -714

$ws.Range("H55").Value = 229.10527
$ws.Range("I55").Value = 201.75
$ws.Range("J55").Value = 276
$ws.Range("K55").Value = 201.75
$ws.Range("L55").Value = 276
$ws.Range("M55").Value = -28.75
$ws.Range("N55").Value = -622

$ws.Range("H61").Value = 1200.2
$ws.Range("I61").Value = 1000.3333
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 1000.3333
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -798.3333
$ws.Range("N61").Value = -1904

$ws.Range("H68").Value = 9908.583000000001
$ws.Range("I68").Value = 26125
$ws.Range("J68").Value = 1800.375
$ws.Range("K68").Value = 26125
$ws.Range("L68").Value = 1800.375
$ws.Range("M68").Value = -25376
$ws.Range("N68").Value = -3298.375

$ws.Range("H71").Value = 9908.583000000001
$ws.Range("I71").Value = 26125
$ws.Range("J71").Value = 1800.375
$ws.Range("K71").Value = 130625
$ws.Range("L71").Value = 9001.875
$ws.Range("M71").Value = -126881
$ws.Range("N71").Value = -16489.875

$ws.Range("H113").Value = 1200.2
$ws.Range("I113").Value = 1000.3333
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1000.3333
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1169.6667
$ws.Range("N113").Value = -5840
